$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.914637720525094
$ws.Range("D2").Value = 7.912942445752905
$ws.Range("E2").Value = 13.14134152479387
$ws.Range("F2").Value = 38.89455146160521
$ws.Range("G2").Value = 45.22471263314228
$ws.Range("H2").Value = 18.09119979066332
$ws.Range("J2").Value = 10.1262676216303
$ws.Range("K2").Value = 18.51904293839889
$ws.Range("N2").Value = 18.09965646904098
$ws.Range("C3").Value = 4.751909275302902
$ws.Range("D3").Value = 7.885275352881312
$ws.Range("E3").Value = 13.09793038102105
$ws.Range("F3").Value = 38.82990988504337
$ws.Range("G3").Value = 45.04017689995115
$ws.Range("H3").Value = 18.12522903229561
$ws.Range("J3").Value = 10.13507144662447
$ws.Range("K3").Value = 18.07559223717352
$ws.Range("N3").Value = 18.17064634286712
$ws.Range("C4").Value = 4.650741037524922
$ws.Range("D4").Value = 7.86924006798673
$ws.Range("E4").Value = 13.07394328694522
$ws.Range("F4").Value = 38.80229116503017
$ws.Range("G4").Value = 44.94382183208309
$ws.Range("H4").Value = 18.15054427988741
$ws.Range("J4").Value = 10.14228136800983
$ws.Range("K4").Value = 17.80187467045279
$ws.Range("N4").Value = 18.21623039557335
$ws.Range("C5").Value = 4.609272908125976
$ws.Range("D5").Value = 7.86294914221526
$ws.Range("E5").Value = 13.06484535954503
$ws.Range("F5").Value = 38.79407176859957
$ws.Range("G5").Value = 44.90883677420206
$ws.Range("H5").Value = 18.16196783370077
$ws.Range("J5").Value = 10.14567280870017
$ws.Range("K5").Value = 17.69014058998834
$ws.Range("N5").Value = 18.23530998413035
$ws.Range("C6").Value = 4.60237483953521
$ws.Range("D6").Value = 7.861919365724727
$ws.Range("E6").Value = 13.06337571310413
$ws.Range("F6").Value = 38.79289021454581
$ws.Range("G6").Value = 44.90328644317236
$ws.Range("H6").Value = 18.16393145626756
$ws.Range("J6").Value = 10.14626332202989
$ws.Range("K6").Value = 17.67158070583631
$ws.Range("N6").Value = 18.23850861186716
$ws.Range("C7").Value = 4.650182659469445
$ws.Range("D7").Value = 7.869154234856659
$ws.Range("E7").Value = 13.07381784023313
$ws.Range("F7").Value = 38.80216802752668
$ws.Range("G7").Value = 44.94333266191063
$ws.Range("H7").Value = 18.15069386424156
$ws.Range("J7").Value = 10.14232527120623
$ws.Range("K7").Value = 17.80036833453956
$ws.Range("N7").Value = 18.21648566768591
$ws.Range("C8").Value = 4.858831426776708
$ws.Range("D8").Value = 7.903207673248414
$ws.Range("E8").Value = 13.12582367672213
$ws.Range("F8").Value = 38.86975856143363
$ws.Range("G8").Value = 45.15757963440393
$ws.Range("H8").Value = 18.10201293392737
$ws.Range("J8").Value = 10.12892849757793
$ws.Range("K8").Value = 18.36653906478169
$ws.Range("N8").Value = 18.1237206032998
$ws.Range("C9").Value = 5.255092439530973
$ws.Range("D9").Value = 7.977355230868094
$ws.Range("E9").Value = 13.24866687459411
$ws.Range("F9").Value = 39.09794308494839
$ws.Range("G9").Value = 45.71107788472106
$ws.Range("H9").Value = 18.04181219566996
$ws.Range("J9").Value = 10.11698712295761
$ws.Range("K9").Value = 19.45829814968292
$ws.Range("N9").Value = 17.95756259141567
$ws.Range("C10").Value = 5.534793147820904
$ws.Range("D10").Value = 8.036059772814584
$ws.Range("E10").Value = 13.35118266172513
$ws.Range("F10").Value = 39.32349712871724
$ws.Range("G10").Value = 46.197032538785
$ws.Range("H10").Value = 18.0193111195871
$ws.Range("J10").Value = 10.11696180851819
$ws.Range("K10").Value = 20.24019641207874
$ws.Range("N10").Value = 17.84497298019231
$ws.Range("C11").Value = 5.658912867591104
$ws.Range("D11").Value = 8.063623988087187
$ws.Range("E11").Value = 13.40037190907282
$ws.Range("F11").Value = 39.43852811574425
$ws.Range("G11").Value = 46.43475692294935
$ws.Range("H11").Value = 18.01383533135715
$ws.Range("J11").Value = 10.11884983375654
$ws.Range("K11").Value = 20.58980123911294
$ws.Range("N11").Value = 17.79578805207119
$ws.Range("C12").Value = 5.705419481936253
$ws.Range("D12").Value = 8.074180011999948
$ws.Range("E12").Value = 13.41935555324273
$ws.Range("F12").Value = 39.48385586041784
$ws.Range("G12").Value = 46.52711540988822
$ws.Range("H12").Value = 18.01244909779535
$ws.Range("J12").Value = 10.11983765010056
$ws.Range("K12").Value = 20.72118312110181
$ws.Range("N12").Value = 17.77745343241186
$ws.Range("C13").Value = 5.695426144131155
$ws.Range("D13").Value = 8.071901422136554
$ws.Range("E13").Value = 13.41525139377867
$ws.Range("F13").Value = 39.47401541902772
$ws.Range("G13").Value = 46.50712139088222
$ws.Range("H13").Value = 18.01271703796294
$ws.Range("J13").Value = 10.11961277671466
$ws.Range("K13").Value = 20.69293439672498
$ws.Range("N13").Value = 17.78138922030506
$ws.Range("C14").Value = 5.662749160003263
$ws.Range("D14").Value = 8.064490104611219
$ws.Range("E14").Value = 13.40192661775544
$ws.Range("F14").Value = 39.44222194170322
$ws.Range("G14").Value = 46.44230882091296
$ws.Range("H14").Value = 18.01370749518539
$ws.Range("J14").Value = 10.11892563561211
$ws.Range("K14").Value = 20.60063097916178
$ws.Range("N14").Value = 17.79427383628862
$ws.Range("C15").Value = 5.642667829780258
$ws.Range("D15").Value = 8.059965666806841
$ws.Range("E15").Value = 13.3938109430113
$ws.Range("F15").Value = 39.42297715510551
$ws.Range("G15").Value = 46.40291181929101
$ws.Range("H15").Value = 18.01440376662605
$ws.Range("J15").Value = 10.11854026548274
$ws.Range("K15").Value = 20.54395776128884
$ws.Range("N15").Value = 17.80220383760948
$ws.Range("C16").Value = 5.526614802503761
$ws.Range("D16").Value = 8.034275209606552
$ws.Range("E16").Value = 13.34801853822917
$ws.Range("F16").Value = 39.31622808329558
$ws.Range("G16").Value = 46.18182709436318
$ws.Range("H16").Value = 18.01976501156364
$ws.Range("J16").Value = 10.11687662381356
$ws.Range("K16").Value = 20.21721505635589
$ws.Range("N16").Value = 17.84822809177997
$ws.Range("C17").Value = 5.454586777997713
$ws.Range("D17").Value = 8.018731150315901
$ws.Range("E17").Value = 13.32057327052257
$ws.Range("F17").Value = 39.25391109083663
$ws.Range("G17").Value = 46.05042639224381
$ws.Range("H17").Value = 18.0242752736622
$ws.Range("J17").Value = 10.11634239153832
$ws.Range("K17").Value = 20.01511424670651
$ws.Range("N17").Value = 17.87698191333041
$ws.Range("C18").Value = 5.41286716579518
$ws.Range("D18").Value = 8.009871832532433
$ws.Range("E18").Value = 13.30502857189963
$ws.Range("F18").Value = 39.21923886306441
$ws.Range("G18").Value = 45.9764202970504
$ws.Range("H18").Value = 18.02731735241346
$ws.Range("J18").Value = 10.11621393558697
$ws.Range("K18").Value = 19.8983056880923
$ws.Range("N18").Value = 17.89371177105831
$ws.Range("C19").Value = 5.398693146527951
$ws.Range("D19").Value = 8.006886333379638
$ws.Range("E19").Value = 13.29980712430595
$ws.Range("F19").Value = 39.20770105698225
$ws.Range("G19").Value = 45.95163480120353
$ws.Range("H19").Value = 18.02842418138199
$ws.Range("J19").Value = 10.11620116160111
$ws.Range("K19").Value = 19.85866334210809
$ws.Range("N19").Value = 17.8994091404793
$ws.Range("C20").Value = 5.462284754863303
$ws.Range("D20").Value = 8.020377475970227
$ws.Range("E20").Value = 13.32346998853699
$ws.Range("F20").Value = 39.26042378770493
$ws.Range("G20").Value = 46.06425191885764
$ws.Range("H20").Value = 18.02374877147001
$ws.Range("J20").Value = 10.11638075591562
$ws.Range("K20").Value = 20.03668773215578
$ws.Range("N20").Value = 17.87390122078703
$ws.Range("C21").Value = 5.672360957066889
$ws.Range("D21").Value = 8.066663828882627
$ws.Range("E21").Value = 13.40583083255582
$ws.Range("F21").Value = 39.45151263853495
$ws.Range("G21").Value = 46.46128290091738
$ws.Range("H21").Value = 18.01339789937125
$ws.Range("J21").Value = 10.11912006323719
$ws.Range("K21").Value = 20.62777102325771
$ws.Range("N21").Value = 17.79048143725486
$ws.Range("C22").Value = 5.806755654592219
$ws.Range("D22").Value = 8.097600087297746
$ws.Range("E22").Value = 13.46173248798899
$ws.Range("F22").Value = 39.58669463924371
$ws.Range("G22").Value = 46.73435884258183
$ws.Range("H22").Value = 18.01064031262935
$ws.Range("J22").Value = 10.12250057135777
$ws.Range("K22").Value = 21.00816143351914
$ws.Range("N22").Value = 17.73765513755459
$ws.Range("C23").Value = 5.735306280777866
$ws.Range("D23").Value = 8.081027961576231
$ws.Range("E23").Value = 13.4317105784165
$ws.Range("F23").Value = 39.51361052762708
$ws.Range("G23").Value = 46.5873900510048
$ws.Range("H23").Value = 18.01174458697509
$ws.Range("J23").Value = 10.12055095921403
$ws.Range("K23").Value = 20.80572195228962
$ws.Range("N23").Value = 17.76569512561295
$ws.Range("C24").Value = 5.458805460746182
$ws.Range("D24").Value = 8.019632931368454
$ws.Range("E24").Value = 13.32215965282236
$ws.Range("F24").Value = 39.25747579943157
$ws.Range("G24").Value = 46.05799660341713
$ws.Range("H24").Value = 18.02398540448201
$ws.Range("J24").Value = 10.11636285479523
$ws.Range("K24").Value = 20.02693627458472
$ws.Range("N24").Value = 17.87529338217018
$ws.Range("C25").Value = 5.149664822497439
$ws.Range("D25").Value = 7.956533910917879
$ws.Range("E25").Value = 13.21324545198774
$ws.Range("F25").Value = 39.02600215954847
$ws.Range("G25").Value = 45.54724398654066
$ws.Range("H25").Value = 18.05429771905179
$ws.Range("J25").Value = 10.11868169046499
$ws.Range("K25").Value = 19.16588604837041
$ws.Range("N25").Value = 18.00083836519124
